# Update the plan and schedule on the "Source Data" sheet.
#
# The iteration table in columns A:D (rows 2-5 and 9-13) lists
# Week / Iteration / Start date / End date. The edit shifts the
# second block of the table (rows 9-13) up by one row into the
# previously-empty row 6, and appends a brand-new "Week 12" row
# at the bottom (row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source Data")

# Row 6 currently only has E6 populated; fill A6:D6 with row 9's
# current (pre-edit) content - this both copies the values and
# their cell formatting/styles.
$ws.Range("A9:D9").Copy($ws.Range("A6:D6"))

# Row 9: Week/Start/End shift up from row 10's current content,
# but the Iteration label is updated to "Iteration 1".
$ws.Range("A10").Copy($ws.Range("A9"))
$ws.Range("C10").Copy($ws.Range("C9"))
$ws.Range("D10").Copy($ws.Range("D9"))
$ws.Range("B3").Copy($ws.Range("B9"))

# Row 10: Week/Start/End shift up from row 11's current content.
$ws.Range("A11").Copy($ws.Range("A10"))
$ws.Range("C11").Copy($ws.Range("C10"))
$ws.Range("D11").Copy($ws.Range("D10"))

# Row 11: Week/Start/End shift up from row 12's current content.
$ws.Range("A12").Copy($ws.Range("A11"))
$ws.Range("C12").Copy($ws.Range("C11"))
$ws.Range("D12").Copy($ws.Range("D11"))

# Row 12: Week/Start/End shift up from row 13's current content.
$ws.Range("A13").Copy($ws.Range("A12"))
$ws.Range("C13").Copy($ws.Range("C12"))
$ws.Range("D13").Copy($ws.Range("D12"))

# Row 13: brand-new "Week 12" entry added to the schedule.
$ws.Range("A13").Value = "Week 12"
$ws.Range("C13").Value = "Monday,15/10/2018"
$ws.Range("D13").Value = "Sunday,22/10/2018"

# Update the selected cell to reflect where the author left off.
$ws.Activate()
$ws.Range("C4").Select()
